$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '''60.622.74'
$ws.Range("E2").Value = '  +2.65%  '

$ws.Range("D3").Formula = '''2.628.87'
$ws.Range("E3").Value = '  +2.36%  '

$ws.Range("D4").Formula = '''0.999'
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").Formula = '''576.00'
$ws.Range("E5").Value = '  +3.99%  '

$ws.Range("D6").Formula = '''142.93'
$ws.Range("E6").Value = '  +1.27%  '

$ws.Range("D7").Formula = '''0.997'
$ws.Range("E7").Value = '  -0.26%  '

$ws.Range("E8").Value = '  +0.62%  '

$ws.Range("D9").Formula = '''2.629.56'
$ws.Range("E9").Value = '  +2.24%  '

$ws.Range("D10").Formula = '''6.50'
$ws.Range("E10").Value = '  -2.10%  '

$ws.Range("E11").Value = '  +2.52%  '

$ws.Range("E12").Value = '  -4.85%  '

$ws.Range("D13").Formula = '''0.367'
$ws.Range("E13").Value = '  +4.68%  '

$ws.Range("D14").Formula = '''3.088.56'
$ws.Range("E14").Value = '  +2.25%  '

$ws.Range("D15").Formula = '''60.634.26'
$ws.Range("E15").Value = '  +2.52%  '

$ws.Range("D16").Formula = '''23.28'
$ws.Range("E16").Value = '  +1.31%  '

$ws.Range("E17").Value = '  +4.76%  '

$ws.Range("D18").Formula = '''2.619.61'
$ws.Range("E18").Value = '  +1.76%  '

$ws.Range("D19").Formula = '''11.27'
$ws.Range("E19").Value = '  +9.63%  '

$ws.Range("E20").Value = '  +2.38%  '

$ws.Range("D21").Formula = '''348.88'
$ws.Range("E21").Value = '  +4.06%  '

$ws.Range("D22").Formula = '''6.91'
$ws.Range("E22").Value = '  +8.39%  '

$ws.Range("D23").Formula = '''0.998'
$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("D24").Formula = '''0.519'
$ws.Range("E24").Value = '  +9.88%  '

$ws.Range("D25").Formula = '''63.25'
$ws.Range("E25").Value = '  +0.86%  '

$ws.Range("D26").Formula = '''0.997'
$ws.Range("E26").Value = '  -0.27%  '

$ws.Range("E27").Value = '  +1.82%  '

$ws.Range("E28").Value = '  +5.67%  '

$ws.Range("D29").Formula = '''0.0₃0795'
$ws.Range("E29").Value = '  +3.68%  '

$ws.Range("D30").Formula = '''1.84'
$ws.Range("E30").Value = '  +10.96%  '

$ws.Range("E31").Value = '  +3.58%  '

$ws.Range("D32").Formula = '''0.997'
$ws.Range("E32").Value = '  -0.12%  '

$ws.Range("D33").Formula = '''161.83'
$ws.Range("E33").Value = '  +2.77%  '

$ws.Range("D34").Formula = '''19.56'
$ws.Range("E34").Value = '  +3.19%  '

$ws.Range("E35").Value = '  +5.29%  '

$ws.Range("D36").Formula = '''0.972'
$ws.Range("E36").Value = '  +9.19%  '

$ws.Range("E37").Value = '  +7.46%  '

$ws.Range("E38").Value = '  +8.76%  '

$ws.Range("E39").Value = '  +1.80%  '

$ws.Range("D40").Formula = '''3.89'
$ws.Range("E40").Value = '  +6.47%  '

$ws.Range("D41").Formula = '''0.846'
$ws.Range("E41").Value = '  -0.57%  '

$ws.Range("D42").Formula = '''301.04'
$ws.Range("E42").Value = '  +3.69%  '

$ws.Range("D43").Formula = '''134.54'
$ws.Range("E43").Value = '  -0.47%  '

$ws.Range("D44").Formula = '''0.997'
$ws.Range("E44").Value = '  -0.19%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Formula = '''19.92'
$ws.Range("E45").Value = '  +5.92%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Formula = '''0.0985'
$ws.Range("E46").Value = '  +1.28%  '

$ws.Range("D47").Formula = '''0.607'
$ws.Range("E47").Value = '  +3.09%  '

$ws.Range("D48").Formula = '''5.01'
$ws.Range("E48").Value = '  +11.32%  '

$ws.Range("E49").Value = '  +3.57%  '

$ws.Range("E50").Value = '  +4.01%  '

$ws.Range("D51").Formula = '''19.95'
$ws.Range("E51").Value = '  +7.61%  '
